# Refresh the live crypto snapshot (price + 1h volume change) for each
# coin row, matching the latest pull from coinranking.com. A few rows
# also swap which coin occupies them (ranking reshuffled this run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds "Price" as plain text (e.g. "1.591.74"); several new
# values look like ordinary decimals (e.g. "211.18"), so COM would
# otherwise auto-coerce them to numbers. Force text via a temporary
# "@" number format, then clear the format again so the cell ends up
# styled exactly like its untouched neighbours.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = "25.992.72"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("D3").Value = "1.591.77"

$ws.Range("E4").Value = "  -0.02%  "

Set-TextValue $ws.Range("D5") "211.18"
$ws.Range("E5").Value = "  +0.47%  "

$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("E7").Value = "  +0.24%  "

Set-TextValue $ws.Range("D8") "0.247"
$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("E9").Value = "  -0.17%  "

Set-TextValue $ws.Range("D10") "18.05"
$ws.Range("E10").Value = "  +0.08%  "

$ws.Range("E11").Value = "  +2.33%  "

$ws.Range("D12").Value = "1.812.07"
$ws.Range("E12").Value = "  +0.45%  "

$ws.Range("D13").Value = "1.595.09"
$ws.Range("E13").Value = "  +0.65%  "

$ws.Range("E14").Value = "  -0.50%  "

$ws.Range("E15").Value = "  +1.17%  "

$ws.Range("D16").Value = "25.975.15"
$ws.Range("E16").Value = "  +0.38%  "

Set-TextValue $ws.Range("D17") "60.59"
$ws.Range("E17").Value = "  +1.34%  "

$ws.Range("E18").Value = "  +0.01%  "

$ws.Range("E19").Value = "  +0.01%  "

Set-TextValue $ws.Range("D20") "202.28"
$ws.Range("E20").Value = "  +5.48%  "

Set-TextValue $ws.Range("D21") "4.25"
$ws.Range("E21").Value = "  +1.60%  "

Set-TextValue $ws.Range("D22") "9.22"
$ws.Range("E22").Value = "  -1.32%  "

Set-TextValue $ws.Range("D23") "6.00"
$ws.Range("E23").Value = "  +1.10%  "

$ws.Range("E24").Value = "  +13.76%  "

Set-TextValue $ws.Range("D25") "143.46"
$ws.Range("E25").Value = "  +1.39%  "

$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("E27").Value = "  -7.48%  "

Set-TextValue $ws.Range("D28") "15.14"
$ws.Range("E28").Value = "  +0.64%  "

Set-TextValue $ws.Range("D29") "6.49"
$ws.Range("E29").Value = "  +0.75%  "

$ws.Range("E30").Value = "  +0.58%  "

$ws.Range("E31").Value = "  +0.91%  "

Set-TextValue $ws.Range("D32") "3.12"
$ws.Range("E32").Value = "  +0.05%  "

Set-TextValue $ws.Range("D33") "2.90"
$ws.Range("E33").Value = "  -3.94%  "

$ws.Range("E34").Value = "  -1.14%  "

Set-TextValue $ws.Range("D35") "2.35"
$ws.Range("E35").Value = "  -0.43%  "

$ws.Range("D36").Value = "1.128.65"
$ws.Range("E36").Value = "  +3.18%  "

$ws.Range("E37").Value = "  +6.83%  "

$ws.Range("E38").Value = "  +0.06%  "

Set-TextValue $ws.Range("D39") "0.793"
$ws.Range("E39").Value = "  +2.37%  "

$ws.Range("E40").Value = "  -1.19%  "

Set-TextValue $ws.Range("D41") "0.489"
$ws.Range("E41").Value = "  -2.23%  "

Set-TextValue $ws.Range("D42") "0.779"
$ws.Range("E42").Value = "  -3.41%  "

$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("D44").Value = "1.724.91"
$ws.Range("E44").Value = "  +0.43%  "

Set-TextValue $ws.Range("D45") "92.22"
$ws.Range("E45").Value = "  -1.18%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D46") "53.69"
$ws.Range("E46").Value = "  +1.18%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D47") "1.48"
$ws.Range("E47").Value = "  -1.04%  "

Set-TextValue $ws.Range("D48") "0.0505"
$ws.Range("E48").Value = "  -0.85%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D49") "0.407"
$ws.Range("E49").Value = "  -0.15%  "

$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
Set-TextValue $ws.Range("D50") "1.01"
$ws.Range("E50").Value = "  +0.41%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₇0947"
$ws.Range("E51").Value = "  -16.08%  "
